$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = "Dox"
$ws.Range("B36").Value = 20
$ws.Range("B37").Select()
